$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 10 appended to the table.
# Column A is an explicit empty text value (use quote-prefix trick so Excel
# stores an actual empty string instead of clearing the cell).
$ws.Range("A10").Value2 = "'"

$ws.Range("B10").Value2 = "يامن "

# Column C looks numeric ("2323") but must be stored as text like the rest
# of the column, so force a text number format before assigning it.
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value2 = "2323"

$ws.Range("D10").Value2 = "الصمود"
$ws.Range("E10").Value2 = "الرحلة 2"
$ws.Range("F10").Value2 = "C3"
$ws.Range("G10").Value2 = "NRC"
$ws.Range("H10").Value2 = "٠١‏/٠٥‏/٢٠٢٥ ٠٦:٥٤:١٧ م"
